{"js": "// The document's \"app build / version\" blurb is being bumped:\n//   \"The app's build is 2.7 (138), where \"2.7 (138)\" stands for major\n//    version 2, minor version 7, revision 0 and build 138.\"\n// becomes:\n//   \"The app's build is 2.7.1 (139), where \"2.7.1 (139)\" stands for major\n//    version 2, minor version 7, revision 1 and build 139.\"\n//\n// i.e. a revision digit (\".1\") is appended to the \"2.7\" version number\n// (both places it appears) and the build number goes from 138 -> 139\n// (both places), and the spelled-out \"revision 0\" becomes \"revision 1\".\n// Word also re-anchors its \"last edit\" (_GoBack) bookmark to sit right\n// after the freshly-typed \"139\", which we mirror below.\n\nasync function replaceAll(body, searchText, replaceText, options) {\n  const results = body.search(searchText, Object.assign({ matchCase: true }, options || {}));\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// \"2.7 (138)\" -> \"2.7.1 (139)\" (occurs twice: the build mention and the\n// quoted repeat of it just after).\nawait replaceAll(body, \"2.7 (138)\", \"2.7.1 (139)\");\n\n// \"revision 0\" -> \"revision 1\"\nawait replaceAll(body, \"revision 0\", \"revision 1\");\n\n// \"build 138\" -> \"build 139\"\nawait replaceAll(body, \"build 138\", \"build 139\");\n\n// Move Word's \"_GoBack\" bookmark (last-edit-location marker) so it sits\n// immediately after the newly typed \"139\".\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // ignore if it didn't exist\n}\n\nconst buildResults = body.search(\"build 139\", { matchCase: true });\nbuildResults.load(\"items\");\nawait context.sync();\nif (buildResults.items.length > 0) {\n  const endRange = buildResults.items[0].getRange(Word.RangeLocation.end);\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The document's \"app build / version\" blurb is being bumped:\n#   \"The app's build is 2.7 (138), where \"2.7 (138)\" stands for major\n#    version 2, minor version 7, revision 0 and build 138.\"\n# becomes:\n#   \"The app's build is 2.7.1 (139), where \"2.7.1 (139)\" stands for major\n#    version 2, minor version 7, revision 1 and build 139.\"\n#\n# i.e. a revision digit (\".1\") is appended to the \"2.7\" version number\n# (both places it appears) and the build number goes from 138 -> 139\n# (both places), and the spelled-out \"revision 0\" becomes \"revision 1\".\n# Word also re-anchors its \"last edit\" (_GoBack) bookmark to sit right\n# after the freshly-typed \"139\", which we mirror below.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($doc, $findText, $replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# \"2.7 (138)\" -> \"2.7.1 (139)\" (occurs twice: the build mention and the\n# quoted repeat of it just after).\nReplace-All $d \"2.7 (138)\" \"2.7.1 (139)\"\n\n# \"revision 0\" -> \"revision 1\"\nReplace-All $d \"revision 0\" \"revision 1\"\n\n# \"build 138\" -> \"build 139\"\nReplace-All $d \"build 138\" \"build 139\"\n\n# Move Word's \"_GoBack\" bookmark (last-edit-location marker) so it sits\n# immediately after the newly typed \"139\".\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"build 139\"\n$found = $find.Execute()\nif ($found) {\n    $matchRange = $find.Parent\n    $endRange = $d.Range($matchRange.End, $matchRange.End)\n    $d.Bookmarks.Add(\"_GoBack\", $endRange) | Out-Null\n}\n"}
